# Generate Report for Archive
# Updates the status of two files (3faa9818... and 51c16a26...) from
# "Ready for handoff" to "In Translation" across the Overview sheet
# (both the zh-cn and de-de status columns) as well as the per-language
# detail sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
# Row 3 -> 3faa9818-6d1c-425a-9aa7-f223f53234ff.md
$overview.Range("B3").Value = "In Translation"
$overview.Range("C3").Value = "In Translation"
# Row 4 -> 51c16a26-eb62-49cb-9756-689eb535fec2.md
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 -> 3faa9818-6d1c-425a-9aa7-f223f53234ff.md (Status column C)
$zhcn.Range("C3").Value = "In Translation"
# Row 4 -> 51c16a26-eb62-49cb-9756-689eb535fec2.md (Status column C)
$zhcn.Range("C4").Value = "In Translation"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
# Row 3 -> 3faa9818-6d1c-425a-9aa7-f223f53234ff.md (Status column C)
$dede.Range("C3").Value = "In Translation"
# Row 4 -> 51c16a26-eb62-49cb-9756-689eb535fec2.md (Status column C)
$dede.Range("C4").Value = "In Translation"
